$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new J-column notes first, in the order that matches the
# target shared-string table layout (234..237).
$ws.Range("J14").Value = "Check colors with Jack"
$ws.Range("J15").Value = "Check with Jack"
$ws.Range("J11").Value = "Check size with Jack"
$ws.Range("J18").Value = "Shield definitions"

$ws.Range("J23").Value = "Check with Jack"
$ws.Range("J25").Value = "Check with Jack"
$ws.Range("J29").Value = "Check with Jack"

# Mark the B/C "Yes" columns for the newly confirmed views.
$ws.Range("B11").Value = "Yes"
$ws.Range("C11").Value = "Yes"

$ws.Range("B14").Value = "Yes"
$ws.Range("C14").Value = "Yes"

$ws.Range("B15").Value = "Yes"
$ws.Range("C15").Value = "Yes"

$ws.Range("B23").Value = "Yes"
$ws.Range("C23").Value = "Yes"

$ws.Range("B25").Value = "Yes"
$ws.Range("C25").Value = "Yes"

$ws.Range("B28").Value = "Yes"
$ws.Range("C28").Value = "Yes"

# Shrink the Dashboards row now that its long note has been replaced.
$ws.Rows.Item(18).RowHeight = 17

$ws.Range("B31").Select()
